# edit.ps1
# Applies the changes described by the commit:
#   "1 meetplayer add DMCActivity 2 foundation add build_ffplay.txt"
# i.e. two new bug-tracker rows are appended to the "bug" sheet (rows 58 & 59),
# a review comment (with adb logcat excerpt) is attached to F59, and the
# active-sheet / view state of the workbook is updated accordingly.

$wb = $excel.ActiveWorkbook

$bug = $wb.Worksheets.Item("bug")
$newfeature = $wb.Worksheets.Item("newfeature")

# ---------------------------------------------------------------------------
# 1) New row 58 (bug #57):
#    id=N/A, platform=android, date=2016-07-25,
#    desc = 乌兰托娅 - 套马杆.ape timestamp discontinuity bug,
#    status=fixed, root cause = audio_decode ret mismatch
# ---------------------------------------------------------------------------
$bug.Range("C58").Value = "N/A"
$bug.Range("D58").Value = "android"
$bug.Range("E58").Value = 20160725
$bug.Range("F58").Value = "乌兰托娅 - 套马杆.ape 时间戳不连续。有1s多间隔"
$bug.Range("G58").Value = "fixed"
$bug.Range("H58").Value = "audio_decode ret =0,0,0,0,xxxx,0,0,0,yyyy 即使output>1，ret也不>0 "

# ---------------------------------------------------------------------------
# 2) New row 59 (bug #58):
#    id=N/A, platform=android, date=2016-08-04,
#    desc = foundation "gotye" build xoplayer aac channel incorrect,
#    status=fixed, root cause = enable-parser / enable-decoder workaround
# ---------------------------------------------------------------------------
$bug.Range("C59").Value = "N/A"
$bug.Range("D59").Value = "android"
$bug.Range("E59").Value = 20160804
$bug.Range("F59").Value = "foundation ""gotye"" build xoplayer aac channel incorrect"
$bug.Range("G59").Value = "fixed"
$bug.Range("H59").Value = "enable-parser=aac,aac_latm enable-decoder=aac_latm"

# ---------------------------------------------------------------------------
# 3) Attach the reviewer comment (logcat excerpt) to F59
# ---------------------------------------------------------------------------
$commentText = "作者:`n08-04 16:22:14.278 7107-7217/? I/FFExtractor: audio codec: codec_id 86018, channels 2, channel_layout 3, sample_rate 44100, sample_fmt -1`n08-04 16:25:36.518 7645-7674/? I/FFExtractor: audio codec: codec_id 86018, channels 1, channel_layout 4, sample_rate 44100, sample_fmt 8"
$comment = $bug.Range("F59").AddComment($commentText)

# ---------------------------------------------------------------------------
# 4) Update view state of "newfeature": it is no longer the active tab,
#    and its viewport scrolls so row 13 is at the top.
# ---------------------------------------------------------------------------
$newfeature.Activate() | Out-Null
$newfeature.Range("F14").Select() | Out-Null
$win = $excel.ActiveWindow
$win.ScrollRow = 13
$win.ScrollColumn = 1

# ---------------------------------------------------------------------------
# 5) Update view state of "bug": it becomes the active tab, scrolled so that
#    row 44 / column D is the top-left visible cell, with H59 selected.
# ---------------------------------------------------------------------------
$bug.Activate() | Out-Null
$bug.Range("H59").Select() | Out-Null
$win = $excel.ActiveWindow
$win.ScrollRow = 44
$win.ScrollColumn = 4

# "bug" (the first sheet) ends up as the active sheet of the workbook.
$bug.Activate() | Out-Null
